# Settings update upon reclustering
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# project_name: dc10 -> iacs_panel_1
$ws.Range("B3").Value = "iacs_panel_1"

# anchor_ids: Anchor -> HC-04
$ws.Range("B9").Value = "HC-04"

# grouping_columns: group, paired_0, ... -> hc_vs_pc_pre_treat, responder_pre_treat, pre_post, pre_post_w_hc
$ws.Range("B15").Value = "hc_vs_pc_pre_treat, responder_pre_treat, pre_post, pre_post_w_hc"

# grouping_orders: DC_0h, DC_0_5h, ... -> HC, S1.1; yes, no; S1.1, S1.2; HC, S1.1, S1.2
$ws.Range("B16").Value = "HC, S1.1; yes, no; S1.1, S1.2; HC, S1.1, S1.2"
# the shorter text now wraps to fewer lines, so the row shrinks accordingly
$ws.Rows.Item(16).RowHeight = 57.6

# data_subsets: DC -> B
$ws.Range("B18").Value = "B"

# clustering_k: 20 -> 30
$ws.Range("B28").Value = 30

# ccp_delta_cutoff: 0.025 -> 0.005
$ws.Range("B30").Value = 0.005

# umap_n: 20 -> 15
$ws.Range("B32").Value = 15

# umap_min_dist: 0.15 -> 0.1
$ws.Range("B33").Value = 0.1

# Update the sheet selection to reflect where the user left off after reclustering
$ws.Range("B30").Select()
